# Fetch Resume API: update Resume sheet with a full entry for the existing
# user (Paramjot Singh) and add a new entry for Aman Kumar, then mark a
# JobApplications row as REJECTED. Finally select the Resume sheet.

function Set-TextValue($cell, $text) {
    # Force a literal text value even when the text looks like a date or a
    # plain number (otherwise Excel auto-converts it to a date serial /
    # numeric value instead of keeping it as text).
    $escaped = $text.Replace('"', '""')
    $cell.Formula = "=""" + $escaped + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.Application.CutCopyMode = $false
}

$wb = $excel.ActiveWorkbook
$resume = $wb.Worksheets.Item("Resume")
$jobApplications = $wb.Worksheets.Item("JobApplications")

# --- Resume sheet: fill in row 2 (Paramjot Singh) with full resume details ---
$resume.Cells.Item(2, 1).Value = "U#00001"
$resume.Cells.Item(2, 2).Value = "Paramjot Singh"
Set-TextValue $resume.Cells.Item(2, 3) "1997-01-26"
$resume.Cells.Item(2, 4).Value = "paramjotsingh966@gmail.com"
$resume.Cells.Item(2, 5).Value = "Reva University"
$resume.Cells.Item(2, 6).Value = "MCA"
Set-TextValue $resume.Cells.Item(2, 7) "2024"
$resume.Cells.Item(2, 8).Value = "Xploria"
$resume.Cells.Item(2, 9).Value = "SDE Internship"
$resume.Cells.Item(2, 10).Value = "6 Months"
$resume.Cells.Item(2, 11).Value = "Frontend`nHTML`nCSS`nJavaScript`nBackend`nJava`nApache Tomcat"

# --- Resume sheet: add row 3 (Aman Kumar) ---
$resume.Cells.Item(3, 1).Value = "U#00002"
$resume.Cells.Item(3, 2).Value = "Aman Kumar"
Set-TextValue $resume.Cells.Item(3, 3) "1998-02-14"
$resume.Cells.Item(3, 4).Value = "aman@gmail.com"
$resume.Cells.Item(3, 5).Value = "Reva University"
$resume.Cells.Item(3, 6).Value = "BCA"
Set-TextValue $resume.Cells.Item(3, 7) "2024"
$resume.Cells.Item(3, 8).Value = "Infosys"
$resume.Cells.Item(3, 9).Value = "Intern"
$resume.Cells.Item(3, 10).Value = "3 Months"
$resume.Cells.Item(3, 11).Value = "Frontend`nReact Js`nAngular Js"

# Widen the Responsibilities column to fit the new content.
$resume.Columns.Item(11).ColumnWidth = 54.33333333333334

# Multi-line responsibilities text auto-expands the row height; restore the
# default (non-custom) row height for both rows.
$resume.Rows.Item(2).AutoFit()
$resume.Rows.Item(3).AutoFit()

# --- JobApplications sheet: mark J#00002 / U#00001 application as REJECTED ---
$jobApplications.Cells.Item(4, 3).Value = "REJECTED"

# --- Make Resume the active sheet, with row 3 selected (the newly added row) ---
$resume.Activate()
$resume.Rows.Item(3).Select()
